$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = -42964.35063398977
$ws.Range("B3").Value = -413.3323703176269
$ws.Range("B4").Value = -209.3780042950629
$ws.Range("B5").Value = 41.90699229350223
$ws.Range("B6").Value = -27.10680516552109
$ws.Range("B7").Value = -670.3845534798008
$ws.Range("B8").Value = -1103.954848074912
$ws.Range("B9").Value = -458.453147936992
$ws.Range("B10").Value = -1537.904399888625
$ws.Range("B11").Value = -169.1976980676666
$ws.Range("B12").Value = -2405.102012099785
